$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shp = $s.Shapes.Item("TextBox 23")
$tr = $shp.TextFrame.TextRange

$old = "post("
$new = "raise("

$target = $tr.Characters(1, $old.Length)
if ($target.Text -eq $old) {
    $target.Text = $new
}
